# Auto-generated Excel COM-interop script
# Applies corrected IFRS financial figures (error solve ifrs list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-6 with corrected financial figures ---
# Row 2
$ws.Range("D2").Value = 4855
$ws.Range("E2").Value = 109
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = -158
$ws.Range("I2").Value = -151
$ws.Range("J2").Value = -7
$ws.Range("K2").Value = 5510
$ws.Range("L2").Value = 3620
$ws.Range("M2").Value = 1890
$ws.Range("N2").Value = 1879
$ws.Range("O2").Value = 11
$ws.Range("P2").Value = 190
$ws.Range("Q2").Value = -214
$ws.Range("R2").Value = 43
$ws.Range("S2").Value = 91
$ws.Range("T2").Value = 88
$ws.Range("U2").Value = -303
$ws.Range("V2").Value = 2798
$ws.Range("W2").Value = 2.24
$ws.Range("X2").Value = -3.26
$ws.Range("Y2").Value = -7.71
$ws.Range("Z2").Value = -2.89
$ws.Range("AA2").Value = 191.48
$ws.Range("AB2").Value = 876.63
$ws.Range("AC2").Value = -734
$ws.Range("AD2").Value = -3.51
$ws.Range("AE2").Value = 9765
$ws.Range("AF2").Value = 0.27
$ws.Range("AG2").Value = 46
$ws.Range("AH2").Value = 1.79
$ws.Range("AI2").Value = -5.84
$ws.Range("AJ2").Value = 20588564

# Row 3
$ws.Range("D3").Value = 4961
$ws.Range("E3").Value = 159
$ws.Range("F3").Value = 159
$ws.Range("G3").Value = -151
$ws.Range("H3").Value = -1102
$ws.Range("I3").Value = -1095
$ws.Range("J3").Value = -7
$ws.Range("K3").Value = 3764
$ws.Range("L3").Value = 2904
$ws.Range("M3").Value = 860
$ws.Range("N3").Value = 856
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 223
$ws.Range("Q3").Value = 300
$ws.Range("R3").Value = 220
$ws.Range("S3").Value = -477
$ws.Range("T3").Value = 40
$ws.Range("U3").Value = 261
$ws.Range("V3").Value = 1682
$ws.Range("W3").Value = 3.21
$ws.Range("X3").Value = -22.21
$ws.Range("Y3").Value = -80.06999999999999
$ws.Range("Z3").Value = -23.77
$ws.Range("AA3").Value = 337.61
$ws.Range("AB3").Value = 259.87
$ws.Range("AC3").Value = -5209
$ws.Range("AD3").Value = -0.41
$ws.Range("AE3").Value = 3955
$ws.Range("AF3").Value = 0.55
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 23000000

# Row 4
$ws.Range("D4").Value = 4445
$ws.Range("E4").Value = 154
$ws.Range("F4").Value = 154
$ws.Range("G4").Value = 236
$ws.Range("H4").Value = 145
$ws.Range("I4").Value = 145
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2861
$ws.Range("L4").Value = 1852
$ws.Range("M4").Value = 1009
$ws.Range("N4").Value = 1005
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 223
$ws.Range("Q4").Value = -353
$ws.Range("R4").Value = 880
$ws.Range("S4").Value = -543
$ws.Range("T4").Value = 38
$ws.Range("U4").Value = -391
$ws.Range("V4").Value = 1171
$ws.Range("W4").Value = 3.45
$ws.Range("X4").Value = 3.26
$ws.Range("Y4").Value = 15.57
$ws.Range("Z4").Value = 4.37
$ws.Range("AA4").Value = 183.63
$ws.Range("AB4").Value = 321.69
$ws.Range("AC4").Value = 630
$ws.Range("AD4").Value = 3.8
$ws.Range("AE4").Value = 4642
$ws.Range("AF4").Value = 0.53
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 23000000

# Row 5
$ws.Range("D5").Value = 3914
$ws.Range("E5").Value = 81
$ws.Range("F5").Value = 81
$ws.Range("G5").Value = 51
$ws.Range("H5").Value = 15
$ws.Range("I5").Value = 15
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2787
$ws.Range("L5").Value = 1778
$ws.Range("M5").Value = 1009
$ws.Range("N5").Value = 1006
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 223
$ws.Range("Q5").Value = -12
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 28
$ws.Range("U5").Value = -40
$ws.Range("V5").Value = 1284
$ws.Range("W5").Value = 2.07
$ws.Range("X5").Value = 0.38
$ws.Range("Y5").Value = 1.52
$ws.Range("Z5").Value = 0.53
$ws.Range("AA5").Value = 176.16
$ws.Range("AB5").Value = 331.36
$ws.Range("AC5").Value = 67
$ws.Range("AD5").Value = 27.58
$ws.Range("AE5").Value = 4645
$ws.Range("AF5").Value = 0.41
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 23000000

# Row 6
$ws.Range("D6").Value = 4132
$ws.Range("E6").Value = 171
$ws.Range("F6").Value = 171
$ws.Range("G6").Value = 154
$ws.Range("H6").Value = 98
$ws.Range("I6").Value = 98
$ws.Range("K6").Value = 3048
$ws.Range("L6").Value = 1925
$ws.Range("M6").Value = 1122
$ws.Range("N6").Value = 1119
$ws.Range("P6").Value = 223
$ws.Range("Q6").Value = 215
$ws.Range("R6").Value = -6
$ws.Range("S6").Value = -65
$ws.Range("T6").Value = 31
$ws.Range("U6").Value = 184
$ws.Range("V6").Value = 1229
$ws.Range("W6").Value = 4.13
$ws.Range("X6").Value = 2.36
$ws.Range("Y6").Value = 9.19
$ws.Range("Z6").Value = 3.34
$ws.Range("AA6").Value = 171.53
$ws.Range("AB6").Value = 369.24
$ws.Range("AC6").Value = 425
$ws.Range("AD6").Value = 4.57
$ws.Range("AE6").Value = 5090
$ws.Range("AF6").Value = 0.39
$ws.Range("AG6").Value = 78
$ws.Range("AH6").Value = 4
$ws.Range("AI6").Value = 17.47
$ws.Range("AJ6").Value = 23000000

# --- Rows 7-9: these companies no longer have financial detail figures;
#     clear all data-value cells (D..AJ), keep the A/B/C labels intact ---
$ws.Range("D7:AJ9").ClearContents()
